$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (S_EQUI / NRANGE_IND): FAKE_RULE was blank, now default("X")
$ws.Range("E4").Value = 'default("X")'

# Row 8 (S_EQUI / BRGEW): FAKE_RULE + Description cleared (rule/description removed)
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()

# Row 9 (S_EQUI / GEWEI): FAKE_RULE was blank, now default("0")
$ws.Range("E9").Value = 'default("0")'

# Shared "foreign_key" rule text updated from foreign_key(S_EQUI.EQUNR) to foreign_key('S_EQUI', 'EQUNR')
# This rule is used by both S_IHPA.EQUNR (row 11) and S_TEXTS_EQUI.EQUNR (row 13, before being replaced below)
$ws.Range("E11").Value = "foreign_key('S_EQUI', 'EQUNR')"

# Row 13 (S_TEXTS_EQUI / EQUNR): FAKE_RULE changes to fk_copy()
$ws.Range("E13").Value = "fk_copy()"

# Row heights recalculated by Excel after the content edits above
$ws.Rows.Item(8).RowHeight = 18.75
$ws.Rows.Item(10).RowHeight = 18.75
